# Applies the content edits captured by the commit:
#  - Slide 4 (TextBox 25): add a trailing space after "...ship liquids along the"
#  - Slide 5 (Title 1 / Assumptions): merge the "Proportionality" run-split back
#    into a single run with the combined wording
#  - Slide 10 (Table 7): rename "Brewery" header cell to "Demand Point"
#  - Slide 11 (Content Placeholder 2): "350 constraints" -> "370 constraints"

$p = $ppt.ActivePresentation

# --- Slide 4: "TextBox 25" --------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(9)          # "TextBox 25"
$tr4 = $shp4.TextFrame.TextRange
$full4 = $tr4.Text
$old4 = " unique costs to ship liquids along the"
$new4 = " unique costs to ship liquids along the "
$start4 = $full4.IndexOf($old4) + 1
$span4 = $tr4.Characters($start4, $old4.Length)
$span4.Text = $new4

# --- Slide 5: "Title 1" (Assumptions list) ----------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(13)         # "Title 1" containing the Assumptions bullets
$tr5 = $shp5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(9)         # "Proportionality: ..." paragraph
$ptext5 = $para5.Text
$old5 = "The value of the objective function and the left-hand side of the constraints are proportional to the level of activity associated with commodities transported and distributed"
$relIdx5 = $ptext5.IndexOf($old5)
$absStart5 = $para5.Start + $relIdx5
$span5 = $tr5.Characters($absStart5, $old5.Length)
$span5.Text = $old5

# --- Slide 10: "Table 7" -----------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(3)        # "Table 7"
$tbl10 = $shp10.Table
$cell10 = $tbl10.Cell(2, 1)
$cell10.Shape.TextFrame.TextRange.Text = "Demand Point"

# --- Slide 11: "Content Placeholder 2" --------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(2)        # "Content Placeholder 2"
$tr11 = $shp11.TextFrame.TextRange
$full11 = $tr11.Text
$old11 = " in Python. The sensitivity report contained information regarding the 350 constraints and the 57 variables. Below, we summarize key findings from our review of this sensitivity report."
$new11 = " in Python. The sensitivity report contained information regarding the 370 constraints and the 57 variables. Below, we summarize key findings from our review of this sensitivity report."
$start11 = $full11.IndexOf($old11) + 1
$span11 = $tr11.Characters($start11, $old11.Length)
$span11.Text = $new11
